$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "DKS" (Desktop Computer) master-data group occupies three consecutive
# rows (the code/name/description rows for eng, ara and fra). Locate it by
# its code in column A rather than a hard-coded row number, then delete the
# whole 3-row block; this also shifts every following group (CMR, SCN, PRT)
# up by three rows automatically.
$dks = $ws.Columns("A").Find("DKS")
if ($dks -ne $null) {
    $block = $ws.Range($dks, $ws.Cells.Item($dks.Row + 2, $dks.Column))
    $block.EntireRow.Delete()
}

# The saved workbook shows cell E10 as the active selection.
$ws.Range("E10").Select()

# The saved worksheet also carries explicit page setup (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
